# Insert a new weekly price-report row right before the current row 40
# (shifts old rows 40-270 down to 41-271) and populate it with the new
# record's data, matching the rest of the dataset's template values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 40..270 down to 41..271, leaving row 40 blank for the
# new record.
$ws.Rows("40:40").Insert()

$ws.Range("A40").Value = 3
$ws.Range("B40").Value = "Femacal de La Calera"
$ws.Range("C40").Value = "Coquimbo"
$ws.Range("D40").Value = 44613
$ws.Range("E40").Value = 5
$ws.Range("F40").Value = 100112039
$ws.Range("G40").Value = "Ciboulette"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 120
$ws.Range("K40").Value = 1500
$ws.Range("L40").Value = 1500
$ws.Range("M40").Value = 1500
$ws.Range("N40").Value = "$/docena de atados"
$ws.Range("O40").Value = "Provincia de Quillota"
$ws.Range("P40").Value = 500
$ws.Range("Q40").Value = 3
$ws.Range("R40").Value = "Hortaliza"
